$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.904957115650177
$ws.Range("B1").Value = 1.703638553619385
$ws.Range("C1").Value = 3.830562829971313
$ws.Range("D1").Value = 3.807269096374512
$ws.Range("E1").Value = 1.206142783164978
